$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new contact rows (values pulled in as new shared strings,
# in the same order they appear in the target sharedStrings table).
$ws.Range("A3").Value = "marie"
$ws.Range("B3").Value = "janne"
$ws.Range("C3").Value = "mariejanne@exemple.com"

$ws.Range("A4").Value = "dupont"
$ws.Range("B4").Value = "pierre"
$ws.Range("C4").Value = "dupontpierre@exemple.com"

# The old hyperlink on C2 is no longer present in the target sheet.
$ws.Range("C2").Hyperlinks.Delete()

# Select the whole populated range, replacing the old H6 selection.
$ws.Range("A1:C4").Select()
